$d = $word.ActiveDocument

# Locate the last paragraph of the document body (ends with "...carrinho correto.")
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Move to the very end of that paragraph's text (before the paragraph mark)
$endRange = $d.Range($lastRange.End - 1, $lastRange.End - 1)

# Insert a new paragraph after it
$endRange.InsertParagraphAfter()

# Grab the newly created (now last) paragraph and set its text via two runs
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertAfter("teste")

$afterTeste = $d.Paragraphs.Last.Range
$afterTeste.Collapse(0)
$afterTeste.InsertAfter("masdfasdf")
